# Add a new worksheet named "Alaaeddin" after the existing "Master" sheet,
# make it the active sheet, put a label in A1, and select A2 on it.
# Also move the selection on the Master sheet to A6 (no longer the
# displayed / active tab).

$wb = $excel.ActiveWorkbook
$master = $wb.Worksheets.Item("Master")

# Add the new sheet right after "Master"
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $master)
$newSheet.Name = "Alaaeddin"

$newSheet.Range("A1").Value = "Data added by Alaaeddin branch"

# Set selection on Master sheet to A6 (leaving the tab as not "selected")
$master.Activate()
$master.Range("A6").Select()

# Make the new sheet the active one (this becomes the workbook's active tab)
$newSheet.Activate()
$newSheet.Range("A2").Select()
